# Auto-generated edit script to update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.748.96"
$ws.Range("E2").Value = "  +1.48%  "
$ws.Range("D3").Value = "'3.797.89"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'697.99"
$ws.Range("E5").Value = "  +7.70%  "
$ws.Range("D6").Value = "'172.59"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("D7").Value = "'3.797.47"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("D11").Value = "'7.28"
$ws.Range("E11").Value = "  +4.53%  "
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("E13").Value = "  +7.33%  "
$ws.Range("D14").Value = "'36.18"
$ws.Range("D15").Value = "'4.437.34"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "'3.797.39"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "'70.717.38"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "'17.72"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "'7.18"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  +15.99%  "
$ws.Range("D22").Value = "'478.24"
$ws.Range("E22").Value = "  +1.74%  "
$ws.Range("D23").Value = "'0.710"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'83.81"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "'12.32"
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "'10.41"
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'2.16"
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "'3.948.77"
$ws.Range("E29").Value = "  +0.44%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("E31").Value = "  +14.71%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.54"
$ws.Range("E32").Value = "  +4.91%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'2.29"
$ws.Range("E33").Value = "  +0.52%  "
$ws.Range("D34").Value = "'0.190"
$ws.Range("E34").Value = "  +8.57%  "
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").Value = "'9.23"
$ws.Range("E36").Value = "  +4.02%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("D39").Value = "'3.41"
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("E41").Value = "  +12.29%  "
$ws.Range("B42").Value = "FLOKI"
$ws.Range("C42").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D42").Value = "'0.000328"
$ws.Range("E42").Value = "  +20.58%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.974"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'163.76"
$ws.Range("E46").Value = "  +4.11%  "
$ws.Range("D47").Value = "'48.83"
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("D48").Value = "'44.49"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("D49").Value = "'0.300"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'1.38"
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "'8.57"
$ws.Range("E51").Value = "  +2.11%  "

Write-Host "Updated crypto data for" 93 "cells"
